$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.686.86"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "1.596.73"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.44"
$ws.Range("D5").ClearFormats()

$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0619"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.52%  "

$ws.Range("E9").Value = "  +0.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.47"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.16%  "

$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").Value = "1.820.57"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").Value = "1.603.26"
$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0₃0769"
$ws.Range("E17").Value = "  +5.68%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.645.47"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "209.05"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.06"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.91"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.83%  "

$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0518"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.75%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("D34").Value = "1.281.32"
$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.615"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.60%  "

$ws.Range("E36").Value = "  -0.44%  "

$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.06"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +17.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.823"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.43"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("E42").Value = "  -0.65%  "

$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.46%  "

$ws.Range("D45").Value = "1.733.11"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.02"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.20%  "

$ws.Range("E47").Value = "  -2.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.101"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("E49").Value = "  +0.54%  "

$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("E51").Value = "  -2.32%  "
